# Update the "changing params" JSON text for rows 2-20 on the active sheet:
# replace the "res_width" / "res_length" pair with a single "res_radius" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = '{"adap_period": ["boundary_code", "permeability", "skin", "res_radius", "length_hor_well_bore", "length_half_fracture", "number_fractures"], "test_period": []}'

for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 2).Value = $newValue
}

# Move the active selection to B2 (matches the saved view state).
$ws.Range("B2").Select()
